$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "Complete CreateTransaction() method"
$ws.Range("D13").Select()
